# July 16 - email password update and run inputs update
# Sheet1 "run inputs" table: remove four stale campaign rows and add a new
# one ("deluxe25offp-redes-July4") just above the redesign summer A/B rows.
# The second table on the same sheet (rows 27-42) is untouched; it simply
# shifts up as a natural consequence of the first table shrinking by three
# rows (17 data rows -> 14 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stale rows from bottom to top so earlier row numbers stay valid:
#   row 15: QA | Sub-D | deluxe25offp-redes
#   row 13: QA | Sub-D | cpcb2017
#   row 9:  QA | Sub-D | deluxe25offp
#   row 7:  QA | Sub-D | cpwbunusedbdbj
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()

# Insert the new campaign row right before the "redes-summera" row (now row 12)
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "QA"
$ws.Range("B12").Value = "Sub-D"
$ws.Range("C12").Value = "deluxe25offp-redes-July4"
$ws.Range("D12").Value = "Kit"
$ws.Range("E12").Value = "Chrome"

# Match the author's final selection in the saved workbook
$ws.Range("C10").Select()
